# "added new code in scheduletype code system"
#
# The FHIR IG publisher (Apache POI) regenerated this workbook with:
#   - Metadata!B8  (Date)  updated to the new generation timestamp
#   - Metadata!B21 (Count) bumped from 5 to 6 concepts
#   - Concepts sheet gained a new row (7) for the "Annen" / Annen" code
#
$wb = $excel.ActiveWorkbook

$wsMeta     = $wb.Worksheets.Item("Metadata")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# --- Metadata sheet -------------------------------------------------------

# Date (Property/Value pair on row 8): bump the generation timestamp.
$wsMeta.Range("B8").Value = "2024-11-07T10:16:17+01:00"

# Count (Property/Value pair on row 21): 5 concepts -> 6 concepts.
# Re-apply row 20's formatting first, then write the value with a leading
# apostrophe so the numeric-looking text is kept as a text value (matching
# the original shared-string "5"/"6") instead of being coerced to a number.
$wsMeta.Range("B20").Copy()
$wsMeta.Range("B21").PasteSpecial(-4122)
$wsMeta.Range("B21").Value = "'6"

# --- Concepts sheet --------------------------------------------------------

# New concept row 7: Level=1, Code=Annen", Display=Annen, Definition=(blank).
# Apply row 6's formatting (border/alignment/style) first, then fill values;
# leading apostrophe keeps the numeric-looking "1" as text, consistent with
# the existing Level column values.
$wsConcepts.Range("A6:D6").Copy()
$wsConcepts.Range("A7:D7").PasteSpecial(-4122)
$wsConcepts.Cells.Item(7, 1).Value = "'1"
$wsConcepts.Cells.Item(7, 2).Value = "Annen`""
$wsConcepts.Cells.Item(7, 3).Value = "Annen"

$excel.CutCopyMode = 0
